# Applies the "Update latest output (run 231)" edit to the optimisation
# result workbook:
#   - Sheet "Schedule": refresh row 2 values and append rows 3 and 4.
#   - Sheet "Detailed": correct historical/forecast values & flags for
#     rows 11-49, and append new forecast rows 50-97 for the next day.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Detailed")

# ---- Sheet1 'Schedule' rows 2-4 ----
$ws1.Range("A2:A4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B2:B4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$arr1 = New-Object 'object[,]' 3,6
$arr1[0,0] = 46064.20833333334
$arr1[0,1] = 46064.66666666666
$arr1[0,2] = 11.0
$arr1[0,3] = 41.58
$arr1[0,4] = 1223.94629175
$arr1[0,5] = 29.43593775252525
$arr1[1,0] = 46065.0
$arr1[1,1] = 46065.25
$arr1[1,2] = 6.0
$arr1[1,3] = 22.68
$arr1[1,4] = 1073.33624775
$arr1[1,5] = 47.32523138227513
$arr1[2,0] = 46065.45833333334
$arr1[2,1] = 46065.66666666666
$arr1[2,2] = 5.0
$arr1[2,3] = 18.9
$arr1[2,4] = 896.26688775
$arr1[2,5] = 47.42152845238095
$ws1.Range("A2:F4").Value() = $arr1

# ---- Sheet2 'Detailed' cell edits (rows 11-49) ----
$ws2.Range("E11").Value() = "OFF"
$ws2.Range("B14").Value() = 122.19086
$ws2.Range("B15").Value() = 119.50455
$ws2.Range("B16").Value() = 59.79985
$ws2.Range("C16").Value() = "historical"
$ws2.Range("B17").Value() = 48.3489
$ws2.Range("C17").Value() = "historical"
$ws2.Range("C18").Value() = "historical"
$ws2.Range("B19").Value() = 48.30146
$ws2.Range("C19").Value() = "historical"
$ws2.Range("B20").Value() = 36.07
$ws2.Range("C20").Value() = "historical"
$ws2.Range("B21").Value() = 50.62162
$ws2.Range("C21").Value() = "historical"
$ws2.Range("B22").Value() = 36.07
$ws2.Range("C22").Value() = "historical"
$ws2.Range("B23").Value() = 36.07
$ws2.Range("C23").Value() = "historical"
$ws2.Range("B24").Value() = 36.07
$ws2.Range("C24").Value() = "historical"
$ws2.Range("B25").Value() = 36.07
$ws2.Range("C25").Value() = "historical"
$ws2.Range("C26").Value() = "historical"
$ws2.Range("B27").Value() = 52.88049
$ws2.Range("C27").Value() = "historical"
$ws2.Range("C28").Value() = "historical"
$ws2.Range("B29").Value() = 62.1786
$ws2.Range("C29").Value() = "historical"
$ws2.Range("B30").Value() = 65.34656
$ws2.Range("C30").Value() = "historical"
$ws2.Range("B31").Value() = 64.89
$ws2.Range("C31").Value() = "historical"
$ws2.Range("B32").Value() = 3.41896
$ws2.Range("C32").Value() = "historical"
$ws2.Range("B33").Value() = 36.06
$ws2.Range("C33").Value() = "historical"
$ws2.Range("E33").Value() = "ON"
$ws2.Range("B34").Value() = 374.03935
$ws2.Range("B35").Value() = 426.55821
$ws2.Range("B36").Value() = 455.85744
$ws2.Range("B37").Value() = 538.56705
$ws2.Range("B38").Value() = 528.66658
$ws2.Range("B39").Value() = 12322.06432
$ws2.Range("B40").Value() = 13169.33477
$ws2.Range("B41").Value() = 663.67173
$ws2.Range("B42").Value() = 511.4643
$ws2.Range("B43").Value() = 117.22663
$ws2.Range("B44").Value() = 189.88
$ws2.Range("B45").Value() = 166.99
$ws2.Range("B46").Value() = 120.89
$ws2.Range("B47").Value() = 145.0
$ws2.Range("B48").Value() = 121.47901
$ws2.Range("B49").Value() = 110.38524

# ---- Sheet2 'Detailed' new rows 50-97 ----
$ws2.Range("A50:A97").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("D50:D97").NumberFormat = "YYYY-MM-DD"
$arr2 = New-Object 'object[,]' 48,5
$arr2[0,0] = 46065.0
$arr2[0,1] = 105.79
$arr2[0,2] = "forecast"
$arr2[0,3] = 46065.0
$arr2[0,4] = "ON"
$arr2[1,0] = 46065.02083333334
$arr2[1,1] = 93.78795
$arr2[1,2] = "forecast"
$arr2[1,3] = 46065.0
$arr2[1,4] = "ON"
$arr2[2,0] = 46065.04166666666
$arr2[2,1] = 93.58072
$arr2[2,2] = "forecast"
$arr2[2,3] = 46065.0
$arr2[2,4] = "ON"
$arr2[3,0] = 46065.0625
$arr2[3,1] = 85.65
$arr2[3,2] = "forecast"
$arr2[3,3] = 46065.0
$arr2[3,4] = "ON"
$arr2[4,0] = 46065.08333333334
$arr2[4,1] = 85.65
$arr2[4,2] = "forecast"
$arr2[4,3] = 46065.0
$arr2[4,4] = "ON"
$arr2[5,0] = 46065.10416666666
$arr2[5,1] = 85.65
$arr2[5,2] = "forecast"
$arr2[5,3] = 46065.0
$arr2[5,4] = "ON"
$arr2[6,0] = 46065.125
$arr2[6,1] = 84.79
$arr2[6,2] = "forecast"
$arr2[6,3] = 46065.0
$arr2[6,4] = "ON"
$arr2[7,0] = 46065.14583333334
$arr2[7,1] = 84.79
$arr2[7,2] = "forecast"
$arr2[7,3] = 46065.0
$arr2[7,4] = "ON"
$arr2[8,0] = 46065.16666666666
$arr2[8,1] = 84.79
$arr2[8,2] = "forecast"
$arr2[8,3] = 46065.0
$arr2[8,4] = "ON"
$arr2[9,0] = 46065.1875
$arr2[9,1] = 85.65
$arr2[9,2] = "forecast"
$arr2[9,3] = 46065.0
$arr2[9,4] = "ON"
$arr2[10,0] = 46065.20833333334
$arr2[10,1] = 100.01
$arr2[10,2] = "forecast"
$arr2[10,3] = 46065.0
$arr2[10,4] = "ON"
$arr2[11,0] = 46065.22916666666
$arr2[11,1] = 110.71902
$arr2[11,2] = "forecast"
$arr2[11,3] = 46065.0
$arr2[11,4] = "ON"
$arr2[12,0] = 46065.25
$arr2[12,1] = 139.62243
$arr2[12,2] = "forecast"
$arr2[12,3] = 46065.0
$arr2[12,4] = "OFF"
$arr2[13,0] = 46065.27083333334
$arr2[13,1] = 138.42
$arr2[13,2] = "forecast"
$arr2[13,3] = 46065.0
$arr2[13,4] = "OFF"
$arr2[14,0] = 46065.29166666666
$arr2[14,1] = 149.47386
$arr2[14,2] = "forecast"
$arr2[14,3] = 46065.0
$arr2[14,4] = "OFF"
$arr2[15,0] = 46065.3125
$arr2[15,1] = 138.42
$arr2[15,2] = "forecast"
$arr2[15,3] = 46065.0
$arr2[15,4] = "OFF"
$arr2[16,0] = 46065.33333333334
$arr2[16,1] = 138.42
$arr2[16,2] = "forecast"
$arr2[16,3] = 46065.0
$arr2[16,4] = "OFF"
$arr2[17,0] = 46065.35416666666
$arr2[17,1] = 105.79
$arr2[17,2] = "forecast"
$arr2[17,3] = 46065.0
$arr2[17,4] = "OFF"
$arr2[18,0] = 46065.375
$arr2[18,1] = 84.79
$arr2[18,2] = "forecast"
$arr2[18,3] = 46065.0
$arr2[18,4] = "OFF"
$arr2[19,0] = 46065.39583333334
$arr2[19,1] = 81.04667
$arr2[19,2] = "forecast"
$arr2[19,3] = 46065.0
$arr2[19,4] = "OFF"
$arr2[20,0] = 46065.41666666666
$arr2[20,1] = 150.98715
$arr2[20,2] = "forecast"
$arr2[20,3] = 46065.0
$arr2[20,4] = "OFF"
$arr2[21,0] = 46065.4375
$arr2[21,1] = 124.3
$arr2[21,2] = "forecast"
$arr2[21,3] = 46065.0
$arr2[21,4] = "OFF"
$arr2[22,0] = 46065.45833333334
$arr2[22,1] = 126.3532
$arr2[22,2] = "forecast"
$arr2[22,3] = 46065.0
$arr2[22,4] = "ON"
$arr2[23,0] = 46065.47916666666
$arr2[23,1] = 84.79
$arr2[23,2] = "forecast"
$arr2[23,3] = 46065.0
$arr2[23,4] = "ON"
$arr2[24,0] = 46065.5
$arr2[24,1] = 76.19007
$arr2[24,2] = "forecast"
$arr2[24,3] = 46065.0
$arr2[24,4] = "ON"
$arr2[25,0] = 46065.52083333334
$arr2[25,1] = 77.66954
$arr2[25,2] = "forecast"
$arr2[25,3] = 46065.0
$arr2[25,4] = "ON"
$arr2[26,0] = 46065.54166666666
$arr2[26,1] = 82.51252
$arr2[26,2] = "forecast"
$arr2[26,3] = 46065.0
$arr2[26,4] = "ON"
$arr2[27,0] = 46065.5625
$arr2[27,1] = 84.79
$arr2[27,2] = "forecast"
$arr2[27,3] = 46065.0
$arr2[27,4] = "ON"
$arr2[28,0] = 46065.58333333334
$arr2[28,1] = 84.79
$arr2[28,2] = "forecast"
$arr2[28,3] = 46065.0
$arr2[28,4] = "ON"
$arr2[29,0] = 46065.60416666666
$arr2[29,1] = 85.09839
$arr2[29,2] = "forecast"
$arr2[29,3] = 46065.0
$arr2[29,4] = "ON"
$arr2[30,0] = 46065.625
$arr2[30,1] = 106.69246
$arr2[30,2] = "forecast"
$arr2[30,3] = 46065.0
$arr2[30,4] = "ON"
$arr2[31,0] = 46065.64583333334
$arr2[31,1] = 110.36191
$arr2[31,2] = "forecast"
$arr2[31,3] = 46065.0
$arr2[31,4] = "ON"
$arr2[32,0] = 46065.66666666666
$arr2[32,1] = 106.15202
$arr2[32,2] = "forecast"
$arr2[32,3] = 46065.0
$arr2[32,4] = "OFF"
$arr2[33,0] = 46065.6875
$arr2[33,1] = 46.66287
$arr2[33,2] = "forecast"
$arr2[33,3] = 46065.0
$arr2[33,4] = "OFF"
$arr2[34,0] = 46065.70833333334
$arr2[34,1] = 83.62278
$arr2[34,2] = "forecast"
$arr2[34,3] = 46065.0
$arr2[34,4] = "OFF"
$arr2[35,0] = 46065.72916666666
$arr2[35,1] = 416.0126
$arr2[35,2] = "forecast"
$arr2[35,3] = 46065.0
$arr2[35,4] = "OFF"
$arr2[36,0] = 46065.75
$arr2[36,1] = 299.99
$arr2[36,2] = "forecast"
$arr2[36,3] = 46065.0
$arr2[36,4] = "OFF"
$arr2[37,0] = 46065.77083333334
$arr2[37,1] = 1040.48406
$arr2[37,2] = "forecast"
$arr2[37,3] = 46065.0
$arr2[37,4] = "OFF"
$arr2[38,0] = 46065.79166666666
$arr2[38,1] = 1061.33673
$arr2[38,2] = "forecast"
$arr2[38,3] = 46065.0
$arr2[38,4] = "OFF"
$arr2[39,0] = 46065.8125
$arr2[39,1] = 411.91618
$arr2[39,2] = "forecast"
$arr2[39,3] = 46065.0
$arr2[39,4] = "OFF"
$arr2[40,0] = 46065.83333333334
$arr2[40,1] = 499.96
$arr2[40,2] = "forecast"
$arr2[40,3] = 46065.0
$arr2[40,4] = "OFF"
$arr2[41,0] = 46065.85416666666
$arr2[41,1] = 160.37894
$arr2[41,2] = "forecast"
$arr2[41,3] = 46065.0
$arr2[41,4] = "OFF"
$arr2[42,0] = 46065.875
$arr2[42,1] = 115.0
$arr2[42,2] = "forecast"
$arr2[42,3] = 46065.0
$arr2[42,4] = "OFF"
$arr2[43,0] = 46065.89583333334
$arr2[43,1] = 119.72055
$arr2[43,2] = "forecast"
$arr2[43,3] = 46065.0
$arr2[43,4] = "OFF"
$arr2[44,0] = 46065.91666666666
$arr2[44,1] = 111.33593
$arr2[44,2] = "forecast"
$arr2[44,3] = 46065.0
$arr2[44,4] = "OFF"
$arr2[45,0] = 46065.9375
$arr2[45,1] = 105.79
$arr2[45,2] = "forecast"
$arr2[45,3] = 46065.0
$arr2[45,4] = "OFF"
$arr2[46,0] = 46065.95833333334
$arr2[46,1] = 74.22007
$arr2[46,2] = "forecast"
$arr2[46,3] = 46065.0
$arr2[46,4] = "OFF"
$arr2[47,0] = 46065.97916666666
$arr2[47,1] = 82.59257
$arr2[47,2] = "forecast"
$arr2[47,3] = 46065.0
$arr2[47,4] = "OFF"
$ws2.Range("A50:E97").Value() = $arr2
